# Add sample_name_prefix argument to data loading (#187)
# Adds a new sample row to the "Samples" sheet (PREFIX_newsample), mirroring
# the existing "BAT-xz971" sample row (same collection date / researcher /
# tissue / collection time / animal id), and makes "Samples" the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# New row 17: same metadata as row 2 (BAT-xz971), but a new, prefixed sample name.
# Force column B to Text format first, matching the existing date-as-text cells
# in column B (numFmt "@"), so "2020-11-19" isn't coerced into a date serial.
$ws.Range("B17").NumberFormat = "@"

$ws.Range("A17").Value = "PREFIX_newsample"
$ws.Range("B17").Value = "2020-11-19"
$ws.Range("C17").Value = "Xianfeng Zhang"
$ws.Range("D17").Value = "BAT"
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 971

# Select A18 under the new last row, matching the edited workbook's saved selection.
[void]$ws.Range("A18").Select()

# Make the Samples sheet the active/selected tab (was Animals).
[void]$ws.Activate()
